$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 3,28

$arr[0,0] = 0.6015231667975562
$arr[0,1] = 0.6220000000000001
$arr[0,2] = 0.5914127684875644
$arr[0,3] = 0.5965
$arr[0,4] = 0.63333111688227
$arr[0,5] = 0.656
$arr[0,6] = 0.6165133457027109
$arr[0,7] = 0.6185
$arr[0,8] = 0.5111399538790977
$arr[0,9] = 0.513
$arr[0,10] = 0.5150948030852934
$arr[0,11] = 0.5165
$arr[0,12] = 0.6308840169499776
$arr[0,13] = 0.6609999999999999
$arr[0,14] = 0.6119448440155988
$arr[0,15] = 0.6144999999999999
$arr[0,16] = 0.6032757859987367
$arr[0,17] = 0.619
$arr[0,18] = 0.5978477645633143
$arr[0,19] = 0.601
$arr[0,20] = 0.6270932346735956
$arr[0,21] = 0.649
$arr[0,22] = 0.6111703256818759
$arr[0,23] = 0.6135
$arr[0,24] = 0.632592589492169
$arr[0,25] = 0.662
$arr[0,26] = 0.6150457195910861
$arr[0,27] = 0.6165

$arr[1,0] = 0.6522627909298521
$arr[1,1] = 0.842
$arr[1,2] = 0.5340715536244627
$arr[1,3] = 0.554
$arr[1,4] = 0.6643208308976678
$arr[1,5] = 0.8530000000000001
$arr[1,6] = 0.5499871638906561
$arr[1,7] = 0.573
$arr[1,8] = 0.6745290942310661
$arr[1,9] = 0.969
$arr[1,10] = 0.5178499292574488
$arr[1,11] = 0.5325
$arr[1,12] = 0.6629241030651151
$arr[1,13] = 0.85
$arr[1,14] = 0.5466654128797196
$arr[1,15] = 0.57
$arr[1,16] = 0.6621633841833641
$arr[1,17] = 0.8779999999999999
$arr[1,18] = 0.5328773928810653
$arr[1,19] = 0.5535
$arr[1,20] = 0.6752951959547212
$arr[1,21] = 0.9229999999999998
$arr[1,22] = 0.5347058109081386
$arr[1,23] = 0.5574999999999999
$arr[1,24] = 0.6645414853783831
$arr[1,25] = 0.875
$arr[1,26] = 0.5378070532789861
$arr[1,27] = 0.5599999999999999

$arr[2,0] = 0.5745815448526939
$arr[2,1] = 0.5660000000000001
$arr[2,2] = 0.5998115209111613
$arr[2,3] = 0.598
$arr[2,4] = 0.674261384613661
$arr[2,5] = 0.6919999999999999
$arr[2,6] = 0.6771920063597272
$arr[2,7] = 0.667
$arr[2,8] = 0.5308052570546892
$arr[2,9] = 0.525
$arr[2,10] = 0.5454438940749439
$arr[2,11] = 0.542
$arr[2,12] = 0.6542885348579179
$arr[2,13] = 0.6450000000000001
$arr[2,14] = 0.6858325444875175
$arr[2,15] = 0.6685
$arr[2,16] = 0.5886385471989254
$arr[2,17] = 0.5830000000000001
$arr[2,18] = 0.6142079792501705
$arr[2,19] = 0.6114999999999999
$arr[2,20] = 0.6735301340586225
$arr[2,21] = 0.6970000000000001
$arr[2,22] = 0.6759493284307114
$arr[2,23] = 0.663
$arr[2,24] = 0.647515729959139
$arr[2,25] = 0.64
$arr[2,26] = 0.6765137304929113
$arr[2,27] = 0.6585

$ws.Range("B4:AC6").Value = $arr
$wb.Save()
